$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Subsetor"
$ws.Range("C1").Value = "Segmento"

$ws.Range("C1").Select()
